# Fix "total marks" calculation error on the marksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking") - Right count corrected from 5 to 4, Wrong marking corrected from -1 to -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total") - recomputed totals after the marking fix
$ws.Range("B12").Value = 56
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "54 / 112"
